$d = $word.ActiveDocument

# --- Update the "Priority" column values in the first table -------------
# Row 3 ("Restaurante" / "Alterar Cadastro"): 3 -> 4
# Row 4 ("Cliente de Restaurante" / "Alterar Cadastro"): 4 -> 5
# Row 6 ("Cliente de Restaurante" / "Sair da fila atual"): 5 -> 3
$t = $d.Tables.Item(1)

$t.Cell(3, 3).Range.Text = "4"
$t.Cell(4, 3).Range.Text = "5"
$t.Cell(6, 3).Range.Text = "3"

# --- Move the "_GoBack" bookmark -----------------------------------------
# It currently sits alone in the empty underlined paragraph right after the
# "Usuário Indefinido = ... Restaurante"" paragraph. It needs to move into
# that paragraph, splitting " Restaurante"" into " Resta" | "urante"".
# Locate the second “Restaurante”” (the one that follows “Cliente de”) and
# split it five characters in, right after "Resta".
$rng = $d.Content
$rng.Find.Execute("Cliente de Restaurante”") | Out-Null
$tail = "urante”"
$splitPoint = $rng.End - $tail.Length
$d.Bookmarks.Add("_GoBack", $d.Range($splitPoint, $splitPoint))
